$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

$ws.Range("K2").Value = 2000

$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G10").Select()
